$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# URL: http://ibm.com/... -> http://linuxforhealth.org/...
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/procedure-type"

# Version: 7.0.0 -> 8.0.0
$meta.Range("B3").Value = "8.0.0"

# Date: 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet 2: "Include from Procedure Type C" ---
$codes = $wb.Worksheets.Item("Include from Procedure Type C")

# System URI: http://ibm.com/... -> http://linuxforhealth.org/...
$codes.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/procedure-type"
